# Add a new forecast vintage column "BB" to the worksheet.
# Column BA holds the most recent existing vintage; BB is the new one
# added one quarter later (date 45986 = 2025-11-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header/date cell BB1, formatted like BA1 ---
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# --- Rows 3-18: new cell BBx carries forward the same value as BAx ---
$carryForwardRows = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18)
$carryForwardValues = @{
    3  = -5.109987415979145
    4  = 2.253603114136604
    5  = 3.738382206110891
    6  = 1.165974434765671
    7  = -0.0426719751787874
    8  = 1.529758493743438
    9  = 1.358758534900462
    10 = 1.664905435092301
    11 = 2.145670176886982
    12 = 1.976124254426503
    13 = 0.7060158009350337
    14 = -4.180878843351332
    15 = 1.312484974417294
    16 = 2.386394320099283
    17 = 0.2104414886460626
    18 = -0.3095793941792935
}
foreach ($r in $carryForwardRows) {
    $ws.Cells.Item($r, 54).Value = $carryForwardValues[$r]
}

# --- Rows 19-21: new cell BBx with newly-forecasted values ---
$ws.Cells.Item(19, 54).Value = -0.08656168856399082
$ws.Cells.Item(20, 54).Value = -0.02867614772544824
$ws.Cells.Item(21, 54).Value = 0
